# ---------------------------------------------------------------------------
# Applies three edits to the letter:
#   1. Updates the letter date from September 19, 2025 to September 21, 2025.
#   2. Splits the "3390 Eichers Pl, Santa Clara CA 95051" mailing-address line
#      (just below the recipient's name) into two lines:
#         "3390 Eichers Pl"
#         "Santa Clara, CA 95051"
#   3. Removes the now-superfluous blank "No Spacing" paragraph that used to
#      sit right after "... Board of Directors".
# ---------------------------------------------------------------------------

$d = $word.ActiveDocument

# 1) Update the letter date (only the standalone date paragraph near the top
#    of the letter, not any other date mentioned in the body).
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -match "^September 19, 2025\r?$") {
        $para.Range.Text = "September 21, 2025"
        break
    }
}

# 2) Split the mailing address paragraph (the one right after the recipient
#    name "Yuning Zheng", i.e. the first occurrence of the combined
#    street/city/state/zip line) into a street line and a city/state/zip line.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -match "^3390 Eichers Pl, Santa Clara CA 95051\r?$") {
        $para.Range.Text = "3390 Eichers Pl"

        # Insert a brand new paragraph right after it, then give it the same
        # paragraph/run formatting (autoSpaceDE/autoSpaceDN + Arial 11pt) as
        # its sibling by inserting raw WordOpenXML for the paragraph.
        $d.Paragraphs.Item($i).Range.InsertParagraphAfter() | Out-Null

        $cityPara = $d.Paragraphs.Item($i + 1)
        $cityXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr><w:t xml:space="preserve">Santa Clara, CA 95051</w:t></w:r></w:p>'
        $insertedOk = $false
        try {
            $cityPara.Range.InsertXML($cityXml)
            $insertedOk = $true
        } catch {
            $insertedOk = $false
        }
        if (-not $insertedOk) {
            # Fallback: at minimum get the right text/paragraph formatting in place
            # even if this runtime doesn't support InsertXML.
            $d.Paragraphs.Item($i + 1).Range.Text = "Santa Clara, CA 95051"
        }
        break
    }
}

# 3) Delete the empty "No Spacing" paragraph that directly follows
#    "...Board of Directors".
for ($i = 1; $i -lt $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -match "Board of Directors") {
        $nextPara = $d.Paragraphs.Item($i + 1)
        if ($nextPara.Range.Text.Trim() -eq "") {
            $nextPara.Range.Delete()
        }
        break
    }
}
